$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny floating point precision of existing A85 date value
$ws.Range("A85").Value = 44398.7692034375

# Append new row 86 with latest data retrieval
$ws.Range("A86").Value = 44399.76833309609
$ws.Range("B86").Value = 80619
$ws.Range("C86").Value = 68060
$ws.Range("D86").Value = 3666
$ws.Range("E86").Value = 2225
$ws.Range("F86").Value = 1604
$ws.Range("G86").Value = 21077
$ws.Range("H86").Value = 1583
$ws.Range("I86").Value = 889
$ws.Range("J86").Value = 195

# Apply the same date/time number format used by the rest of column A
$ws.Range("A86").NumberFormat = $ws.Range("A85").NumberFormat
